$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force column D (Price) cells to remain text so numeric-looking
# values such as "232.71" or "0.999" are not auto-converted to numbers,
# matching the original inlineStr text cells.
$dCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D9', 'D10', 'D12', 'D13', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D26', 'D28', 'D29', 'D30', 'D31', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated Coin / Link / Price / Volume(1h) values
$ws.Range('D2').Value = '43.797.68'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '2.349.75'
$ws.Range('E3').Value = '  +3.64%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '232.71'
$ws.Range('E5').Value = '  +1.06%  '
$ws.Range('D6').Value = '0.650'
$ws.Range('E6').Value = '  +3.79%  '
$ws.Range('D7').Value = '66.04'
$ws.Range('E7').Value = '  +4.66%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '0.451'
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').Value = '0.0975'
$ws.Range('E10').Value = '  -4.21%  '
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '26.92'
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('D13').Value = '2.693.30'
$ws.Range('E13').Value = '  +3.31%  '
$ws.Range('E14').Value = '  -1.05%  '
$ws.Range('D15').Value = '15.48'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').Value = '6.17'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('D17').Value = '0.851'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').Value = '2.346.99'
$ws.Range('E18').Value = '  +4.05%  '
$ws.Range('D19').Value = '43.642.56'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').Value = '0.0₃0985'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').Value = '74.09'
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '6.28'
$ws.Range('E22').Value = '  +3.56%  '
$ws.Range('D23').Value = '250.00'
$ws.Range('E23').Value = '  -0.90%  '
$ws.Range('D24').Value = '3.81'
$ws.Range('E24').Value = '  +15.22%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('D28').Value = '9.93'
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('D29').Value = '22.42'
$ws.Range('E29').Value = '  +7.48%  '
$ws.Range('D30').Value = '174.95'
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('D31').Value = '1.45'
$ws.Range('E31').Value = '  +4.26%  '
$ws.Range('E32').Value = '  -4.89%  '
$ws.Range('E33').Value = '  +1.53%  '
$ws.Range('D34').Value = '5.00'
$ws.Range('E34').Value = '  +4.29%  '
$ws.Range('D35').Value = '0.0690'
$ws.Range('E35').Value = '  -1.57%  '
$ws.Range('D36').Value = '4.98'
$ws.Range('E36').Value = '  +1.47%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '2.46'
$ws.Range('E37').Value = '  +6.48%  '
$ws.Range('D38').Value = '3.73'
$ws.Range('E38').Value = '  -2.94%  '
$ws.Range('B39').Value = 'THORChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D39').Value = '6.56'
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('D41').Value = '9.22'
$ws.Range('E41').Value = '  +11.71%  '
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D43').Value = '17.87'
$ws.Range('E43').Value = '  +1.18%  '
$ws.Range('D44').Value = '1.19'
$ws.Range('E44').Value = '  +10.40%  '
$ws.Range('D45').Value = '99.59'
$ws.Range('E45').Value = '  +1.40%  '
$ws.Range('D46').Value = '0.0956'
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('D47').Value = '1.20'
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('D48').Value = '4.33'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('D49').Value = '1.449.56'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').Value = '10.02'
$ws.Range('E50').Value = '  -1.93%  '
$ws.Range('D51').Value = '2.30'
$ws.Range('E51').Value = '  +0.01%  '

Write-Output "Updated cryptos list"
